# Auto-generated edit script applying the scrape update for LÍNEA 141 (17/01/2026)
# Updates timestamps, re-sorts/swaps a handful of tied rows, and appends newly scraped rows.

$wb = $excel.ActiveWorkbook

# ============================================================
# Sheet 1: LP1912
# ============================================================
$ws1 = $wb.Worksheets.Item("LP1912")
$ws1.Range("A2").Value = "Última actualización: 10:48:14"
$ws1.Range("A3").Value = "Total filas: 168"

$ws1.Range("A46").Value = "06:33:46"
$ws1.Range("B46").Value = "07:59"
$ws1.Range("C46").Value = "11_ETCHEVERRY"
$ws1.Range("D46").Value = 86
$ws1.Range("E46").Value = "LP1912"

$ws1.Range("A47").Value = "07:12:53"
$ws1.Range("B47").Value = "07:59"
$ws1.Range("C47").Value = "23_HERNANDEZ"
$ws1.Range("D47").Value = 47
$ws1.Range("E47").Value = "LP1912"

$ws1.Range("A109").Value = "10:04:17"
$ws1.Range("B109").Value = "10:04"
$ws1.Range("C109").Value = "215C_EL PATO"
$ws1.Range("D109").Value = 0
$ws1.Range("E109").Value = "LP1912"

$ws1.Range("A110").Value = "08:46:25"
$ws1.Range("B110").Value = "10:04"
$ws1.Range("C110").Value = "14_ABASTO"
$ws1.Range("D110").Value = 78
$ws1.Range("E110").Value = "LP1912"

$ws1.Range("A130").Value = "10:36:18"
$ws1.Range("B130").Value = "10:37"
$ws1.Range("C130").Value = "16_SANTA ANA"
$ws1.Range("D130").Value = 1
$ws1.Range("E130").Value = "LP1912"

$ws1.Range("A131").Value = "08:39:08"
$ws1.Range("B131").Value = "10:37"
$ws1.Range("C131").Value = "16_P MOR-SANTA ANA"
$ws1.Range("D131").Value = 118
$ws1.Range("E131").Value = "LP1912"

$ws1.Range("A136").Value = "10:48:14"
$ws1.Range("B136").Value = "10:48"
$ws1.Range("C136").Value = "10_OLMOS"
$ws1.Range("D136").Value = 0
$ws1.Range("E136").Value = "LP1912"

$ws1.Range("A137").Value = "10:48:14"
$ws1.Range("B137").Value = "10:49"
$ws1.Range("C137").Value = "16_SANTA ANA"
$ws1.Range("D137").Value = 1
$ws1.Range("E137").Value = "LP1912"

$ws1.Range("A138").Value = "10:04:17"
$ws1.Range("B138").Value = "10:51"
$ws1.Range("C138").Value = "15_ABASTO"
$ws1.Range("D138").Value = 47
$ws1.Range("E138").Value = "LP1912"

$ws1.Range("A139").Value = "10:36:18"
$ws1.Range("B139").Value = "10:54"
$ws1.Range("C139").Value = "10_OLMOS"
$ws1.Range("D139").Value = 18
$ws1.Range("E139").Value = "LP1912"

$ws1.Range("A140").Value = "09:21:49"
$ws1.Range("B140").Value = "10:56"
$ws1.Range("C140").Value = "27_EL RETIRO"
$ws1.Range("D140").Value = 95
$ws1.Range("E140").Value = "LP1912"

$ws1.Range("A141").Value = "10:04:17"
$ws1.Range("B141").Value = "10:57"
$ws1.Range("C141").Value = "27_EL RETIRO"
$ws1.Range("D141").Value = 53
$ws1.Range("E141").Value = "LP1912"

$ws1.Range("A142").Value = "09:21:49"
$ws1.Range("B142").Value = "11:01"
$ws1.Range("C142").Value = "17_ROMERO"
$ws1.Range("D142").Value = 100
$ws1.Range("E142").Value = "LP1912"

$ws1.Range("A143").Value = "10:36:18"
$ws1.Range("B143").Value = "11:03"
$ws1.Range("C143").Value = "23_HERNANDEZ"
$ws1.Range("D143").Value = 27
$ws1.Range("E143").Value = "LP1912"

$ws1.Range("A144").Value = "09:21:49"
$ws1.Range("B144").Value = "11:04"
$ws1.Range("C144").Value = "14_ABASTO"
$ws1.Range("D144").Value = 103
$ws1.Range("E144").Value = "LP1912"

$ws1.Range("A145").Value = "10:04:17"
$ws1.Range("B145").Value = "11:05"
$ws1.Range("C145").Value = "14_ABASTO"
$ws1.Range("D145").Value = 61
$ws1.Range("E145").Value = "LP1912"

$ws1.Range("A146").Value = "10:36:18"
$ws1.Range("B146").Value = "11:11"
$ws1.Range("C146").Value = "15_ABASTO"
$ws1.Range("D146").Value = 35
$ws1.Range("E146").Value = "LP1912"

$ws1.Range("A147").Value = "10:04:17"
$ws1.Range("B147").Value = "11:11"
$ws1.Range("C147").Value = "23_HERNANDEZ"
$ws1.Range("D147").Value = 67
$ws1.Range("E147").Value = "LP1912"

$ws1.Range("A148").Value = "09:21:49"
$ws1.Range("B148").Value = "11:14"
$ws1.Range("C148").Value = "225_C ROCA-H SUR"
$ws1.Range("D148").Value = 113
$ws1.Range("E148").Value = "LP1912"

$ws1.Range("A149").Value = "09:21:49"
$ws1.Range("B149").Value = "11:20"
$ws1.Range("C149").Value = "215C_EL PATO"
$ws1.Range("D149").Value = 119
$ws1.Range("E149").Value = "LP1912"

$ws1.Range("A150").Value = "10:04:17"
$ws1.Range("B150").Value = "11:21"
$ws1.Range("C150").Value = "215C_EL PATO"
$ws1.Range("D150").Value = 77
$ws1.Range("E150").Value = "LP1912"

$ws1.Range("A151").Value = "10:48:14"
$ws1.Range("B151").Value = "11:21"
$ws1.Range("C151").Value = "10_OLMOS"
$ws1.Range("D151").Value = 33
$ws1.Range("E151").Value = "LP1912"

$ws1.Range("A152").Value = "10:36:18"
$ws1.Range("B152").Value = "11:22"
$ws1.Range("C152").Value = "10_OLMOS"
$ws1.Range("D152").Value = 46
$ws1.Range("E152").Value = "LP1912"

$ws1.Range("A153").Value = "10:36:18"
$ws1.Range("B153").Value = "11:24"
$ws1.Range("C153").Value = "11_ETCHEVERRY"
$ws1.Range("D153").Value = 48
$ws1.Range("E153").Value = "LP1912"

$ws1.Range("A154").Value = "10:36:18"
$ws1.Range("B154").Value = "11:25"
$ws1.Range("C154").Value = "16_P MOR-SANTA ANA"
$ws1.Range("D154").Value = 49
$ws1.Range("E154").Value = "LP1912"

$ws1.Range("A155").Value = "10:04:17"
$ws1.Range("B155").Value = "11:25"
$ws1.Range("C155").Value = "11_ETCHEVERRY"
$ws1.Range("D155").Value = 81
$ws1.Range("E155").Value = "LP1912"

$ws1.Range("A156").Value = "10:04:17"
$ws1.Range("B156").Value = "11:30"
$ws1.Range("C156").Value = "15X38_ABASTO"
$ws1.Range("D156").Value = 86
$ws1.Range("E156").Value = "LP1912"

$ws1.Range("A157").Value = "10:48:14"
$ws1.Range("B157").Value = "11:32"
$ws1.Range("C157").Value = "23_HERNANDEZ"
$ws1.Range("D157").Value = 44
$ws1.Range("E157").Value = "LP1912"

$ws1.Range("A158").Value = "10:36:18"
$ws1.Range("B158").Value = "11:33"
$ws1.Range("C158").Value = "23_HERNANDEZ"
$ws1.Range("D158").Value = 57
$ws1.Range("E158").Value = "LP1912"

$ws1.Range("A159").Value = "10:48:14"
$ws1.Range("B159").Value = "11:33"
$ws1.Range("C159").Value = "16_SANTA ANA"
$ws1.Range("D159").Value = 45
$ws1.Range("E159").Value = "LP1912"

$ws1.Range("A160").Value = "10:04:17"
$ws1.Range("B160").Value = "11:34"
$ws1.Range("C160").Value = "10_OLMOS"
$ws1.Range("D160").Value = 90
$ws1.Range("E160").Value = "LP1912"

$ws1.Range("A161").Value = "10:36:18"
$ws1.Range("B161").Value = "11:35"
$ws1.Range("C161").Value = "16_SANTA ANA"
$ws1.Range("D161").Value = 59
$ws1.Range("E161").Value = "LP1912"

$ws1.Range("A162").Value = "10:04:17"
$ws1.Range("B162").Value = "11:37"
$ws1.Range("C162").Value = "16_P MOR-SANTA ANA"
$ws1.Range("D162").Value = 93
$ws1.Range("E162").Value = "LP1912"

$ws1.Range("A163").Value = "10:04:17"
$ws1.Range("B163").Value = "11:40"
$ws1.Range("C163").Value = "215A_EL PATO"
$ws1.Range("D163").Value = 96
$ws1.Range("E163").Value = "LP1912"

$ws1.Range("A164").Value = "10:04:17"
$ws1.Range("B164").Value = "11:45"
$ws1.Range("C164").Value = "215B_EL PATO"
$ws1.Range("D164").Value = 101
$ws1.Range("E164").Value = "LP1912"

$ws1.Range("A165").Value = "10:04:17"
$ws1.Range("B165").Value = "11:54"
$ws1.Range("C165").Value = "225_GOMEZ"
$ws1.Range("D165").Value = 110
$ws1.Range("E165").Value = "LP1912"

$ws1.Range("A166").Value = "10:48:14"
$ws1.Range("B166").Value = "12:07"
$ws1.Range("C166").Value = "14_ABASTO"
$ws1.Range("D166").Value = 79
$ws1.Range("E166").Value = "LP1912"

$ws1.Range("A167").Value = "10:36:18"
$ws1.Range("B167").Value = "12:29"
$ws1.Range("C167").Value = "215C_EL PATO"
$ws1.Range("D167").Value = 113
$ws1.Range("E167").Value = "LP1912"

$ws1.Range("A168").Value = "10:36:18"
$ws1.Range("B168").Value = "12:30"
$ws1.Range("C168").Value = "11_ETCHEVERRY"
$ws1.Range("D168").Value = 114
$ws1.Range("E168").Value = "LP1912"

$ws1.Range("A169").Value = "10:36:18"
$ws1.Range("B169").Value = "12:31"
$ws1.Range("C169").Value = "16_P MOR-SANTA ANA"
$ws1.Range("D169").Value = 115
$ws1.Range("E169").Value = "LP1912"

$ws1.Range("A170").Value = "10:48:14"
$ws1.Range("B170").Value = "12:31"
$ws1.Range("C170").Value = "11_ETCHEVERRY"
$ws1.Range("D170").Value = 103
$ws1.Range("E170").Value = "LP1912"

$ws1.Range("A171").Value = "10:48:14"
$ws1.Range("B171").Value = "12:37"
$ws1.Range("C171").Value = "27_EL RETIRO"
$ws1.Range("D171").Value = 109
$ws1.Range("E171").Value = "LP1912"

$ws1.Range("A172").Value = "10:48:14"
$ws1.Range("B172").Value = "12:40"
$ws1.Range("C172").Value = "15X38_ABASTO"
$ws1.Range("D172").Value = 112
$ws1.Range("E172").Value = "LP1912"

$ws1.Range("A173").Value = "10:48:14"
$ws1.Range("B173").Value = "12:43"
$ws1.Range("C173").Value = "14_ABASTO"
$ws1.Range("D173").Value = 115
$ws1.Range("E173").Value = "LP1912"

# ============================================================
# Sheet 2: LP1912-215 (only the refreshed timestamp changes)
# ============================================================
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws2.Range("A2").Value = "Última actualización: 10:48:14"

# ============================================================
# Sheet 3: 6203-6173
# ============================================================
$ws3 = $wb.Worksheets.Item("6203-6173")
$ws3.Range("A2").Value = "Última actualización: 10:48:14"
$ws3.Range("A3").Value = "Total filas: 20"
$ws3.Range("A25").Value = "10:48:14"
$ws3.Range("B25").Value = "12:44"
$ws3.Range("C25").Value = "215C_LA PLATA"
$ws3.Range("D25").Value = 116
$ws3.Range("E25").Value = "L6203"

